$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create a new "2022-Q4" worksheet right after "总计", by duplicating the
#    existing "2022-Q3" sheet (so that formatting / sheetPr / page margins
#    are preserved faithfully) and overwriting its data with the Q4 figures.
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $summarySheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Row 2 - 招商量化精选股票A
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'001917"
$q4Sheet.Range("C2").Value = "招商量化精选股票A"
$q4Sheet.Range("D2").Value = "'5.91"
$q4Sheet.Range("E2").Value = "'94.08"
$q4Sheet.Range("F2").Value = "'1.44"
$q4Sheet.Range("G2").Value = "'0.0851"
$q4Sheet.Range("H2").Value = 7

# Row 3 - 招商量化精选股票C
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'007950"
$q4Sheet.Range("C3").Value = "招商量化精选股票C"
$q4Sheet.Range("D3").Value = "'5.28"
$q4Sheet.Range("E3").Value = "'94.08"
$q4Sheet.Range("F3").Value = "'1.44"
$q4Sheet.Range("G3").Value = "'0.0760"
$q4Sheet.Range("H3").Value = 7

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new summary row for 2022-Q4 above the
#    existing quarters, pushing everything else down by one row.
# ---------------------------------------------------------------------------
$ws = $summarySheet

# Shift rows 2-5 down to rows 3-6 (go bottom-up so we don't clobber data).
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Range("B$dest").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dest").Value = $ws.Range("C$r").Value2
    $ws.Range("D$dest").Value = $ws.Range("D$r").Value2
}

# New row 6 needs column A formatted/styled like the other index cells; copy
# the style down from row 5 before writing the final index value.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 4

$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# New 2022-Q4 summary row.
$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.16
